# Auto-generated edit script: applies scheduled market-data refresh values
# to the Leve profit sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2229.4119
$ws.Range("I113").Value = 2086.889
$ws.Range("J113").Value = 2389.75
$ws.Range("K113").Value = 2086.889
$ws.Range("L113").Value = 2389.75
$ws.Range("M113").Value = 1167.111
$ws.Range("N113").Value = -8897.75
$ws.Range("H125").Value = 1866.4
$ws.Range("I125").Value = 1032
$ws.Range("J125").Value = 2075
$ws.Range("K125").Value = 9288
$ws.Range("L125").Value = 18675
$ws.Range("M125").Value = -6828
$ws.Range("N125").Value = -23595
$ws.Range("H129").Value = 873.59375
$ws.Range("J129").Value = 1006.88
$ws.Range("L129").Value = 3020.64
$ws.Range("N129").Value = -13020.64
$ws.Range("H138").Value = 3667.8708
$ws.Range("I138").Value = 2177.2812
$ws.Range("J138").Value = 5257.8335
$ws.Range("K138").Value = 6531.8436
$ws.Range("L138").Value = 15773.5005
$ws.Range("M138").Value = -1391.8436
$ws.Range("N138").Value = -26053.5005

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 389793.03
$ws.Range("I32").Value = 2945.7612
$ws.Range("J32").Value = 4092474
$ws.Range("K32").Value = 2945.7612
$ws.Range("L32").Value = 4092474
$ws.Range("M32").Value = -2658.7612
$ws.Range("N32").Value = -4093048
$ws.Range("H122").Value = 3288.1904
$ws.Range("I122").Value = 3605.8823
$ws.Range("J122").Value = 1938
$ws.Range("K122").Value = 10817.6469
$ws.Range("L122").Value = 5814
$ws.Range("M122").Value = -8367.6469
$ws.Range("N122").Value = -10714
$ws.Range("H132").Value = 1173.3903
$ws.Range("I132").Value = 941.2059
$ws.Range("J132").Value = 2301.1428
$ws.Range("K132").Value = 2823.6177
$ws.Range("L132").Value = 6903.428400000001
$ws.Range("M132").Value = -293.6177000000002
$ws.Range("N132").Value = -11963.4284

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8930700
$ws.Range("I31").Value = 12501794
$ws.Range("K31").Value = 12501794
$ws.Range("M31").Value = -12501499
$ws.Range("H34").Value = 8930700
$ws.Range("I34").Value = 12501794
$ws.Range("K34").Value = 12501794
$ws.Range("M34").Value = -12501592
$ws.Range("H58").Value = 1365.5143
$ws.Range("I58").Value = 1194.238
$ws.Range("J58").Value = 1622.4286
$ws.Range("K58").Value = 1194.238
$ws.Range("L58").Value = 1622.4286
$ws.Range("M58").Value = -991.2380000000001
$ws.Range("N58").Value = -2028.4286
$ws.Range("H136").Value = 1365.5143
$ws.Range("I136").Value = 1194.238
$ws.Range("J136").Value = 1622.4286
$ws.Range("K136").Value = 3582.714
$ws.Range("L136").Value = 4867.2858
$ws.Range("M136").Value = -1032.714
$ws.Range("N136").Value = -9967.2858

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 742.63635
$ws.Range("I5").Value = 454.32144
$ws.Range("J5").Value = 1247.1875
$ws.Range("K5").Value = 1362.96432
$ws.Range("L5").Value = 3741.5625
$ws.Range("M5").Value = -1250.96432
$ws.Range("N5").Value = -3965.5625
$ws.Range("H122").Value = 481.9091
$ws.Range("I122").Value = 327.54544
$ws.Range("J122").Value = 636.2727
$ws.Range("K122").Value = 2947.90896
$ws.Range("L122").Value = 5726.454299999999
$ws.Range("M122").Value = -497.9089599999998
$ws.Range("N122").Value = -10626.4543
$ws.Range("H131").Value = 7353779
$ws.Range("I131").Value = 775.7143
$ws.Range("J131").Value = 9260114
$ws.Range("K131").Value = 2327.1429
$ws.Range("L131").Value = 27780342
$ws.Range("M131").Value = 2712.8571
$ws.Range("N131").Value = -27790422
$ws.Range("H135").Value = 742.63635
$ws.Range("I135").Value = 454.32144
$ws.Range("J135").Value = 1247.1875
$ws.Range("K135").Value = 4088.89296
$ws.Range("L135").Value = 11224.6875
$ws.Range("M135").Value = -1553.89296
$ws.Range("N135").Value = -16294.6875
$ws.Range("H136").Value = 1989.1538
$ws.Range("I136").Value = 1084.4445
$ws.Range("J136").Value = 4024.75
$ws.Range("K136").Value = 3253.3335
$ws.Range("L136").Value = 12074.25
$ws.Range("M136").Value = 1846.6665
$ws.Range("N136").Value = -22274.25

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 1098.75
$ws.Range("I19").Value = 697.5
$ws.Range("J19").Value = 1500
$ws.Range("K19").Value = 697.5
$ws.Range("L19").Value = 1500
$ws.Range("M19").Value = -409.5
$ws.Range("N19").Value = -2076
$ws.Range("H64").Value = 28000
$ws.Range("J64").Value = 28000
$ws.Range("L64").Value = 28000
$ws.Range("N64").Value = -28496
$ws.Range("H67").Value = 28000
$ws.Range("J67").Value = 28000
$ws.Range("L67").Value = 28000
$ws.Range("N67").Value = -29716
$ws.Range("H122").Value = 2003.7333
$ws.Range("I122").Value = 1861.2273
$ws.Range("J122").Value = 2395.625
$ws.Range("K122").Value = 5583.6819
$ws.Range("L122").Value = 7186.875
$ws.Range("M122").Value = -3133.6819
$ws.Range("N122").Value = -12086.875
$ws.Range("H126").Value = 8773454
$ws.Range("I126").Value = 1292.2222
$ws.Range("J126").Value = 16668400
$ws.Range("K126").Value = 3876.6666
$ws.Range("L126").Value = 50005200
$ws.Range("M126").Value = -1406.6666
$ws.Range("N126").Value = -50010140

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2310
$ws.Range("I40").Value = 2080
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 2080
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -1944
$ws.Range("N40").Value = -3272
$ws.Range("H61").Value = 1912.6875
$ws.Range("I61").Value = 1536.6364
$ws.Range("K61").Value = 1536.6364
$ws.Range("M61").Value = -1334.6364
$ws.Range("H113").Value = 1912.6875
$ws.Range("I113").Value = 1536.6364
$ws.Range("K113").Value = 1536.6364
$ws.Range("M113").Value = 633.3635999999999
$ws.Range("H132").Value = 1851.0695
$ws.Range("I132").Value = 1861.0189
$ws.Range("J132").Value = 1823.3158
$ws.Range("K132").Value = 5583.0567
$ws.Range("L132").Value = 5469.9474
$ws.Range("M132").Value = -3053.0567
$ws.Range("N132").Value = -10529.9474
$ws.Range("H136").Value = 3157.2954
$ws.Range("I136").Value = 2421.7942
$ws.Range("J136").Value = 5658
$ws.Range("K136").Value = 7265.382599999999
$ws.Range("L136").Value = 16974
$ws.Range("M136").Value = -4715.382599999999
$ws.Range("N136").Value = -22074

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3531.2942
$ws.Range("I81").Value = 1913.2
$ws.Range("J81").Value = 5842.857
$ws.Range("K81").Value = 3826.4
$ws.Range("L81").Value = 11685.714
$ws.Range("M81").Value = -2765.4
$ws.Range("N81").Value = -13807.714
$ws.Range("H84").Value = 3531.2942
$ws.Range("I84").Value = 1913.2
$ws.Range("J84").Value = 5842.857
$ws.Range("K84").Value = 19132
$ws.Range("L84").Value = 58428.57
$ws.Range("M84").Value = -13828
$ws.Range("N84").Value = -69036.57000000001
$ws.Range("H132").Value = 16667562
$ws.Range("I132").Value = 18383204
$ws.Range("J132").Value = 1329.2858
$ws.Range("K132").Value = 55149612
$ws.Range("L132").Value = 3987.8574
$ws.Range("M132").Value = -55147082
$ws.Range("N132").Value = -9047.857400000001
$ws.Range("H136").Value = 1312.8518
$ws.Range("I136").Value = 1641.2858
$ws.Range("J136").Value = 959.1539
$ws.Range("K136").Value = 4923.857400000001
$ws.Range("L136").Value = 2877.4617
$ws.Range("M136").Value = -2373.857400000001
$ws.Range("N136").Value = -7977.4617
